$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D price cells to Text format so that numeric-looking
# strings (e.g. "1.005", "219.02") are stored as text, matching the source data,
# instead of being auto-converted to floating point numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.585.65'
$ws.Range("E2").Value = '  -7.28%  '

$ws.Range("D3").Value = '1.695.53'
$ws.Range("E3").Value = '  -5.86%  '

$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.26%  '

$ws.Range("D5").Value = '219.02'
$ws.Range("E5").Value = '  -5.53%  '

$ws.Range("D6").Value = '0.5086'
$ws.Range("E6").Value = '  -13.84%  '

$ws.Range("D7").Value = '1.005'
$ws.Range("E7").Value = '  +0.18%  '

$ws.Range("D8").Value = '0.2646'
$ws.Range("E8").Value = '  -4.45%  '

$ws.Range("D9").Value = '22.10'
$ws.Range("E9").Value = '  -5.11%  '

$ws.Range("D10").Value = '0.06286'
$ws.Range("E10").Value = '  -7.58%  '

$ws.Range("D11").Value = '0.07367'
$ws.Range("E11").Value = '  -1.98%  '

$ws.Range("D12").Value = '1.696.07'
$ws.Range("E12").Value = '  -8.61%  '

$ws.Range("D13").Value = '4.511'
$ws.Range("E13").Value = '  -5.27%  '

$ws.Range("D14").Value = '0.5833'
$ws.Range("E14").Value = '  -6.20%  '

$ws.Range("D15").Value = '1.925.60'
$ws.Range("E15").Value = '  -5.86%  '

$ws.Range("D16").Value = '0.000008386'
$ws.Range("E16").Value = '  -8.22%  '

$ws.Range("D17").Value = '65.51'
$ws.Range("E17").Value = '  -13.22%  '

$ws.Range("D18").Value = '26.614.48'
$ws.Range("E18").Value = '  -7.06%  '

$ws.Range("D19").Value = '5.011'
$ws.Range("E19").Value = '  -8.16%  '

$ws.Range("E20").Value = '  +0.24%  '

$ws.Range("E21").Value = '  -4.42%  '

$ws.Range("D22").Value = '185.98'
$ws.Range("E22").Value = '  -11.60%  '

$ws.Range("D23").Value = '6.268'
$ws.Range("E23").Value = '  -8.06%  '

$ws.Range("E24").Value = '  +0.22%  '

$ws.Range("D25").Value = '144.54'
$ws.Range("E25").Value = '  -6.02%  '

$ws.Range("D26").Value = '7.505'
$ws.Range("E26").Value = '  -4.31%  '

$ws.Range("D27").Value = '0.1156'
$ws.Range("E27").Value = '  -8.84%  '

$ws.Range("D28").Value = '15.67'
$ws.Range("E28").Value = '  -4.65%  '

$ws.Range("D29").Value = '1.342'
$ws.Range("E29").Value = '  -5.43%  '

$ws.Range("D30").Value = '0.05647'
$ws.Range("E30").Value = '  -8.64%  '

$ws.Range("D31").Value = '1.332'
$ws.Range("E31").Value = '  -6.40%  '

$ws.Range("D33").Value = '3.492'
$ws.Range("E33").Value = '  -6.57%  '

$ws.Range("D34").Value = '1.642'
$ws.Range("E34").Value = '  -4.98%  '

$ws.Range("E35").Value = '  -3.55%  '

$ws.Range("D36").Value = '0.6029'
$ws.Range("E36").Value = '  -6.08%  '

$ws.Range("E37").Value = '  -5.34%  '

$ws.Range("D38").Value = '2.679'
$ws.Range("E38").Value = '  -1.23%  '

$ws.Range("D39").Value = '0.01611'
$ws.Range("E39").Value = '  -4.92%  '

$ws.Range("D40").Value = '1.095.56'
$ws.Range("E40").Value = '  -4.53%  '

$ws.Range("D41").Value = '0.8599'
$ws.Range("E41").Value = '  -2.67%  '

$ws.Range("D42").Value = '5.848'
$ws.Range("E42").Value = '  -10.67%  '

$ws.Range("E43").Value = '  -0.44%  '

$ws.Range("D44").Value = '99.34'
$ws.Range("E44").Value = '  -0.72%  '

$ws.Range("D45").Value = '1.852.43'
$ws.Range("E45").Value = '  -5.23%  '

$ws.Range("D46").Value = '0.00000000110'
$ws.Range("E46").Value = '  -1.27%  '

$ws.Range("D47").Value = '56.71'
$ws.Range("E47").Value = '  -5.95%  '

$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '8.144'
$ws.Range("E48").Value = '  -2.44%  '

$ws.Range("B49").Value = 'Frax'
$ws.Range("C49").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D49").Value = '1.004'
$ws.Range("E49").Value = '  +0.26%  '

$ws.Range("D50").Value = '0.05240'
$ws.Range("E50").Value = '  -4.25%  '

$ws.Range("D51").Value = '0.4328'
$ws.Range("E51").Value = '  -3.30%  '

# Restore the default (Normal) style on the D column so no stray explicit
# per-cell style/number-format reference is left behind on the cells.
$ws.Range("D2:D51").Style = "Normal"
